# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Re-sorts the 92 existing "Periodo Mora" detail rows (16-107) into
#   ascending chronological order (1712 .. 2507) instead of descending.
# - Adds 4 new detail rows (108-111) for period 2508 covering 3 new
#   workers + the existing worker.
# - Updates the summary block (VALOR MORA, Cant. Trabajadores, Cant.
#   Periodos).
# - Moves the signature block down from rows 112/113 to rows 116/117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Re-sort the period column (E16:E107) into ascending order. Every
#    other column in these rows (doc type, doc #, name, value, salary)
#    already holds the correct repeated value and does not change.
# ---------------------------------------------------------------------
$periods = @("1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212","2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312","2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 5).Value = $periods[$i]
}

# ---------------------------------------------------------------------
# 2. Row 107 used to be the last (bold-bordered) row of the table. It
#    is now a normal row, so restyle it like the rows above it, and
#    stash its original "last row" formatting to re-use on row 111,
#    the new last row.
# ---------------------------------------------------------------------
$ws.Range("B107:J107").Copy()
$ws.Range("B111:J111").PasteSpecial(-4122)

$ws.Range("B106:J106").Copy()
$ws.Range("B107:J107").PasteSpecial(-4122)
$ws.Range("B106:J106").Copy()
$ws.Range("B108:J108").PasteSpecial(-4122)
$ws.Range("B106:J106").Copy()
$ws.Range("B109:J109").PasteSpecial(-4122)
$ws.Range("B106:J106").Copy()
$ws.Range("B110:J110").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. Fill in the 4 new detail rows for period 2508.
# ---------------------------------------------------------------------
$ws.Range("B108").Value = "CC"
$ws.Range("C108").Value = "73109672"
$ws.Range("D108").Value = "EDGAR CASTILLO LEDESMA"
$ws.Range("E108").Value = "2508"
$ws.Range("F108").Value = 100000
$ws.Range("G108").Value = 2500000

$ws.Range("B109").Value = "CC"
$ws.Range("C109").Value = "3805234"
$ws.Range("D109").Value = "ALEXANDER LUGO ARROYO"
$ws.Range("E109").Value = "2508"
$ws.Range("F109").Value = 32000
$ws.Range("G109").Value = 800000

$ws.Range("B110").Value = "CC"
$ws.Range("C110").Value = "45764059"
$ws.Range("D110").Value = "VERENA MARIA LUGO ARROYO"
$ws.Range("E110").Value = "2508"
$ws.Range("F110").Value = 80000
$ws.Range("G110").Value = 2000000

$ws.Range("B111").Value = "PPT"
$ws.Range("C111").Value = "898055"
$ws.Range("D111").Value = "MAGDELY SANYAIR ARANGUREN ARRIECHE"
$ws.Range("E111").Value = "2508"
$ws.Range("F111").Value = 80000
$ws.Range("G111").Value = 2000000

# Column D needs to be wide enough for the longest new name.
$ws.Columns("D").ColumnWidth = 39

# ---------------------------------------------------------------------
# 4. Update the summary block above the table.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 3236000   # VALOR MORA
$ws.Range("C13").Value = 4         # Cant. Trabajadores
$ws.Range("F13").Value = 93        # Cant. Periodos

# ---------------------------------------------------------------------
# 5. Move the signature block from rows 112/113 down to rows 116/117.
# ---------------------------------------------------------------------
$ws.Range("B112:C112").UnMerge()
$ws.Range("B112:C112").Clear()
$ws.Range("H112:J112").UnMerge()
$ws.Range("H112:J112").Clear()
$ws.Range("B113:C113").UnMerge()
$ws.Range("B113:C113").Clear()
$ws.Range("H113:J113").UnMerge()
$ws.Range("H113:J113").Clear()

$ws.Range("B116").Value = "___________________________________"
$ws.Range("B116:C116").Merge()
$ws.Range("H116").Value = "___________________________________"
$ws.Range("H116:J116").Merge()

$ws.Range("B117").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("B117:C117").Merge()
$ws.Range("H117").Value = "FIRMA DEL REPRESENTANTE LEGAL"
$ws.Range("H117:J117").Merge()
